# Ajuste na extração de remanescentes concluído
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows 25-29 (Marlos de Medeiros Chaves moves to the top of the
#     block, Venucia Bruna Magalhaes Pereira moves after Claudia) -----------
$rows25to29 = @(
    @("Marlos de Medeiros Chaves", 0, 90, 0, 40, 80, 90, 0, 0, 300, 20),
    @("Adriana Costa Bacelo", 50, 0, 110, 0, 10, 35, 0, 80, 285, 30),
    @("Anna Carolina Machado Marinho", 0, 15, 0, 0, 150, 0, 15, 90, 270, 20),
    @("Claudia Stutz Zubieta", 0, 0, 0, 0, 0, 90, 90, 80, 260, 0),
    @("Venúcia Bruna Magalhães Pereira", 80, 80, 60, 0, 40, 0, 0, 0, 260, 0)
)

$r = 25
foreach ($rowData in $rows25to29) {
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    $r++
}

# --- Swap rows 35 and 36 -----------------------------------------------
$rows35to36 = @(
    @("Giovanny Augusto Camacho Antevere Mazzarotto", 80, 0, 0, 0, 5, 0, 15, 0, 100, 0),
    @("Margareth Borges Coutinho Gallo", 0, 0, 0, 0, 100, 0, 0, 0, 100, 25)
)

$r = 35
foreach ($rowData in $rows35to36) {
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    $r++
}

# --- Append new row 40 ---------------------------------------------------
# Copy the formatting of the preceding "Autor" cell (bold, bordered,
# centered/top aligned) onto the new row's name cell before setting values.
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)  # xlPasteFormats

$newRow = @("Fernanda Savicki de Almeida", 0, 0, 0, 0, 0, 0, 5, 0, 5, 0)
for ($c = 0; $c -lt $newRow.Length; $c++) {
    $ws.Cells.Item(40, $c + 1).Value = $newRow[$c]
}
